$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph (paragraph 1). Built with the same shape as the other
#    body paragraphs in this document (leading empty run, then a bold
#    label run, then a plain run with the rest of the sentence).
# -----------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null
$metaPara = $d.Paragraphs.Item(2)

$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Book of Anunnaki, a 10-payline slot game with a free spins bonus round and expanding special symbol. Play for free.</w:t></w:r></w:p>"
$metaPara.Range.InsertXML($metaXml) | Out-Null

# -----------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the
#    document ("Play Book of Anunnaki Free: A Mystical Egyptian-Themed
#    Slot Game").
# -----------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Font.Bold = 1
$find.Text = "Play Book of Anunnaki Free: A Mystical Egyptian-Themed Slot Game"
$found = $find.Execute()
if ($found) {
    $dupRange = $find.Parent
    $dupRange.Expand(4) | Out-Null
    $dupRange.Delete()
}

# -----------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new
#    image-prompt text, keeping its existing formatting. Scope the
#    search to the last paragraph only so the similarly-worded text
#    inside the new "Meta description" paragraph is left untouched.
# -----------------------------------------------------------------
$oldText = "Read our review of Book of Anunnaki, a 10-payline slot game with a free spins bonus round and expanding special symbol. Play for free."
$newText = "Create a feature image that will catch the attention of players of Book of Anunnaki. The image should be in cartoon style and have a happy Maya warrior wearing glasses. The warrior should be standing in front of an ancient Egyptian temple, holding the Book of Anunnaki in one hand while smiling at the rewards in the other hand. The background should be a desert scene with pyramids and camels. Make sure to include the game logo and name in the image to make it recognizable to players. The image should be colorful and playful to attract potential players to try the game."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldText, $true, $false, $false, $false, $false,
                              $true, 1, $false, $newText, 2) | Out-Null
